# The uploaded workbook replaces the row for client 2128 with a new
# client 2506 (row 19 of the "Planilha1" sheet): the EMPRD code in
# column A and the EMPREENDIMENTO description in column B are updated.
# Columns C (ADM) and D (UF) are left untouched, matching the source
# diff (row 19 keeps ADM="Maria Eduarda" / UF="RJ").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

$ws.Range("A19").Value = 2506
$ws.Range("B19").Value = "2506 - KATIA FERRIRA DE BARROS"
